$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-10-16 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-17 Thursday", 2)

$d.Content.Find.Execute("748÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "495÷4=", 2)
$d.Content.Find.Execute("604÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "735÷4=", 2)
$d.Content.Find.Execute("331÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "572÷3=", 2)
$d.Content.Find.Execute("996÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "162÷8=", 2)
$d.Content.Find.Execute("379÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "747÷3=", 2)

$d.Content.Find.Execute("338÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "433÷3=", 2)
$d.Content.Find.Execute("846÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "858÷6=", 2)
$d.Content.Find.Execute("956÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "643÷8=", 2)
$d.Content.Find.Execute("540÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "676÷2=", 2)
$d.Content.Find.Execute("352÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "783÷6=", 2)

$d.Content.Find.Execute("187÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "274÷6=", 2)
$d.Content.Find.Execute("276÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "496÷9=", 2)
$d.Content.Find.Execute("798÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "744÷2=", 2)
$d.Content.Find.Execute("142÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "644÷2=", 2)
$d.Content.Find.Execute("568÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "810÷6=", 2)

$d.Content.Find.Execute("292÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "149÷7=", 2)
$d.Content.Find.Execute("786÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "889÷9=", 2)
$d.Content.Find.Execute("293÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "348÷4=", 2)
$d.Content.Find.Execute("575÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "655÷6=", 2)
$d.Content.Find.Execute("925÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "774÷4=", 2)

$d.Content.Find.Execute("655÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "674÷9=", 2)
$d.Content.Find.Execute("231÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "803÷4=", 2)
$d.Content.Find.Execute("662÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "875÷9=", 2)
$d.Content.Find.Execute("506÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "345÷8=", 2)
$d.Content.Find.Execute("170÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "482÷8=", 2)
